# variantSummary.xlsx update:
#  - remove 3 rows for variants that came from duplicate genomes
#    (POS 2334, 18404, 18974 - each was a singleton nSubjects=1/nSamples=1 row)
#  - adjust nSubjects (col B) / nSamples (col C) counts for the positions
#    whose supporting-sample counts changed once the duplicate genomes were
#    collapsed out
#  - re-sort the table by nSamples (desc) then POS (asc), same ordering rule
#    the report uses throughout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Delete the rows for the POS values that no longer appear at all.
#    Find each by scanning column A, then delete from the bottom up so
#    earlier row numbers stay valid while we work.
# ---------------------------------------------------------------------
$positionsToRemove = @(2334, 18404, 18974)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$rowsToDelete = @()
for ($r = 2; $r -le $lastRow; $r++) {
    $posVal = $ws.Cells.Item($r, 1).Value2
    if ($positionsToRemove -contains [int]$posVal) {
        $rowsToDelete += $r
    }
}

$rowsToDelete = $rowsToDelete | Sort-Object -Descending
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# ---------------------------------------------------------------------
# 2) Update nSubjects (B) / nSamples (C) for the positions whose counts
#    changed because duplicate genomes were collapsed.
# ---------------------------------------------------------------------
$updates = @{
    241   = @(29, 53)
    3037  = @(29, 53)
    8782  = @(29, 53)
    14408 = @(29, 53)
    18060 = @(29, 53)
    23403 = @(29, 53)
    28144 = @(29, 53)
    25563 = @(27, 48)
    1059  = @(23, 39)
    18255 = @(5, 9)
    18877 = @(3, 8)
    7386  = @(1, 4)
    11083 = @(2, 3)
    23994 = @(1, 4)
    1500  = @(1, 2)
    3090  = @(1, 2)
    13501 = @(1, 2)
    5570  = @(1, 1)
    18512 = @(1, 1)
}

for ($r = 2; $r -le $lastRow; $r++) {
    $posVal = [int]$ws.Cells.Item($r, 1).Value2
    if ($updates.ContainsKey($posVal)) {
        $pair = $updates[$posVal]
        $ws.Cells.Item($r, 2).Value = $pair[0]
        $ws.Cells.Item($r, 3).Value = $pair[1]
    }
}

# ---------------------------------------------------------------------
# 3) Re-sort the data (A2:E<lastRow>) by nSamples (C) descending, then
#    POS (A) ascending - matching the report's ordering convention.
# ---------------------------------------------------------------------
$sortRange = $ws.Range("A2:E$lastRow")
$keyC = $ws.Range("C2:C$lastRow")
$keyA = $ws.Range("A2:A$lastRow")

$sortRange.Sort($keyC, 2, $keyA, [Type]::Missing, 1, [Type]::Missing, 1, 1)

Write-Host "Done. lastRow=$lastRow"
